$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 68

# --- Text columns: Date, Time, Weekday, Week -------------------------------
# In the source file these are plain text (e.g. "2024-01-17", "02" with a
# leading zero) rather than real Excel dates/numbers. Assigning such
# look-alike strings straight to .Value makes Excel auto-convert them to a
# date serial / number, so instead we build each one as a formula that
# *evaluates* to the literal text, then copy/paste-special the computed
# values back onto themselves. That bakes in the plain string without
# triggering the "smart" re-parsing (and without leaving a stray
# quote-prefix style behind).
$ws.Cells.Item($row, 1).Formula = '="2024-01-17"'
$ws.Cells.Item($row, 2).Formula = '="17:23:26"'
$ws.Cells.Item($row, 3).Formula = '="Wednesday"'
$ws.Cells.Item($row, 4).Formula = '="02"'

$textRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 4))
$textRange.Copy()
$textRange.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# --- Numeric columns: Beijing .. Wuhan --------------------------------------
$ws.Cells.Item($row, 5).Value = 138985
$ws.Cells.Item($row, 6).Value = 139585
$ws.Cells.Item($row, 7).Value = 171112
$ws.Cells.Item($row, 8).Value = 148626
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 119552
$ws.Cells.Item($row, 11).Value = 222761
$ws.Cells.Item($row, 12).Value = 254816
$ws.Cells.Item($row, 13).Value = 184861
$ws.Cells.Item($row, 14).Value = 110302
$ws.Cells.Item($row, 15).Value = 41310
$ws.Cells.Item($row, 16).Value = 30928
$ws.Cells.Item($row, 17).Value = 73447
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42563
$ws.Cells.Item($row, 20).Value = -1
